$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new log row documenting a "last_update" change for the bevnat-info document.
$row = 54
$ws.Cells.Item($row, 1).Value = 1759778989
$ws.Cells.Item($row, 2).Value = "update"
$ws.Cells.Item($row, 3).Value = "doc"
$ws.Cells.Item($row, 4).Value = "bevnat-info"
$ws.Cells.Item($row, 6).Value = "last_update"

# Force these numeric-looking values to be stored as text (matching the
# source data, which keeps all log values as strings) without leaving a
# residual number-format style on the cell.
$ws.Cells.Item($row, 7).Value = "'1706219962"
$ws.Cells.Item($row, 7).ClearFormats()
$ws.Cells.Item($row, 8).Value = "'1706239962"
$ws.Cells.Item($row, 8).ClearFormats()
